$wb = $excel.ActiveWorkbook

# --- Update scraped_at timestamps on the "snapshot" sheet (K2:K31) ---
$ws = $wb.Worksheets.Item("snapshot")
$ws.Range("K2").Value = "2025-11-30T07:01:42.649296+00:00"
$ws.Range("K3").Value = "2025-11-30T07:01:44.930065+00:00"
$ws.Range("K4").Value = "2025-11-30T07:01:44.930098+00:00"
$ws.Range("K5").Value = "2025-11-30T07:01:48.228919+00:00"
$ws.Range("K6").Value = "2025-11-30T07:01:50.482953+00:00"
$ws.Range("K7").Value = "2025-11-30T07:01:53.151136+00:00"
$ws.Range("K8").Value = "2025-11-30T07:01:55.471089+00:00"
$ws.Range("K9").Value = "2025-11-30T07:02:01.119922+00:00"
$ws.Range("K10").Value = "2025-11-30T07:02:01.119950+00:00"
$ws.Range("K11").Value = "2025-11-30T07:02:03.365314+00:00"
$ws.Range("K12").Value = "2025-11-30T07:02:06.132647+00:00"
$ws.Range("K13").Value = "2025-11-30T07:02:06.132680+00:00"
$ws.Range("K14").Value = "2025-11-30T07:02:08.729639+00:00"
$ws.Range("K15").Value = "2025-11-30T07:02:11.036570+00:00"
$ws.Range("K16").Value = "2025-11-30T07:02:11.036601+00:00"
$ws.Range("K17").Value = "2025-11-30T07:02:11.036621+00:00"
$ws.Range("K18").Value = "2025-11-30T07:02:13.227284+00:00"
$ws.Range("K19").Value = "2025-11-30T07:02:13.227321+00:00"
$ws.Range("K20").Value = "2025-11-30T07:02:13.227345+00:00"
$ws.Range("K21").Value = "2025-11-30T07:02:15.518413+00:00"
$ws.Range("K22").Value = "2025-11-30T07:02:15.518445+00:00"
$ws.Range("K23").Value = "2025-11-30T07:02:17.945083+00:00"
$ws.Range("K24").Value = "2025-11-30T07:02:17.945117+00:00"
$ws.Range("K25").Value = "2025-11-30T07:02:17.945142+00:00"
$ws.Range("K26").Value = "2025-11-30T07:02:20.253176+00:00"
$ws.Range("K27").Value = "2025-11-30T07:02:25.528438+00:00"
$ws.Range("K28").Value = "2025-11-30T07:02:30.137711+00:00"
$ws.Range("K29").Value = "2025-11-30T07:02:30.137743+00:00"
$ws.Range("K30").Value = "2025-11-30T07:02:32.912727+00:00"
$ws.Range("K31").Value = "2025-11-30T07:02:32.912754+00:00"

# --- Remove the now-stale INJURED_NEW row from "new_injured" (row 2) ---
$ws2 = $wb.Worksheets.Item("new_injured")
$ws2.Rows.Item(2).Delete()
